$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers formatted as localized-looking
# strings (inline/shared strings in the source data), several of which
# parse as plain decimals (e.g. "577.58"). Excel auto-converts a numeric
# looking string typed into a General-formatted cell into a real number,
# so force those particular cells to Text format first to keep them as
# literal strings, matching the source data's text-cell representation.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    # Reset the visual style back to Normal (xf 0) now that the literal
    # text is committed, so no stray number-format style lingers on the
    # cell (source cells carry no explicit style in columns B:E).
    $c.Style = "Normal"
}

Set-TextValue "D2" "63.636.97"

Set-TextValue "D3" "3.437.85"
$ws.Range("E3").Value = "  +7.54%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue "D5" "577.58"
$ws.Range("E5").Value = "  +7.33%  "

Set-TextValue "D6" "156.85"
$ws.Range("E6").Value = "  +7.60%  "

$ws.Range("E7").Value = "  -0.05%  "

Set-TextValue "D8" "3.448.42"
$ws.Range("E8").Value = "  +7.67%  "

Set-TextValue "D9" "0.534"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("E10").Value = "  +2.98%  "

$ws.Range("E11").Value = "  +8.68%  "

$ws.Range("E12").Value = "  +0.80%  "

Set-TextValue "D13" "4.028.71"
$ws.Range("E13").Value = "  +7.46%  "

$ws.Range("E14").Value = "  -0.47%  "

Set-TextValue "D15" "0.0000187"

Set-TextValue "D16" "27.26"
$ws.Range("E16").Value = "  +5.71%  "

Set-TextValue "D17" "63.791.62"
$ws.Range("E17").Value = "  +6.24%  "

Set-TextValue "D18" "3.437.78"
$ws.Range("E18").Value = "  +7.57%  "

Set-TextValue "D19" "6.43"
$ws.Range("E19").Value = "  +2.56%  "

$ws.Range("E20").Value = "  +7.29%  "

Set-TextValue "D21" "8.50"
$ws.Range("E21").Value = "  +3.60%  "

Set-TextValue "D22" "391.95"
$ws.Range("E22").Value = "  +5.77%  "

$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("E24").Value = "  +3.23%  "

Set-TextValue "D25" "72.10"
$ws.Range("E25").Value = "  +3.59%  "

Set-TextValue "D26" "0.0000107"
$ws.Range("E26").Value = "  +22.45%  "

Set-TextValue "D27" "9.53"
$ws.Range("E27").Value = "  +10.77%  "

$ws.Range("E28").Value = "  +7.16%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +7.81%  "

Set-TextValue "D31" "6.60"
$ws.Range("E31").Value = "  +7.77%  "

Set-TextValue "D32" "1.36"
$ws.Range("E32").Value = "  +14.95%  "

Set-TextValue "D33" "5.74"
$ws.Range("E33").Value = "  +8.88%  "

Set-TextValue "D34" "23.48"
$ws.Range("E34").Value = "  +4.49%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("E37").Value = "  +9.41%  "

Set-TextValue "D38" "158.71"
$ws.Range("E38").Value = "  -0.20%  "

Set-TextValue "D39" "28.16"
$ws.Range("E39").Value = "  +5.89%  "

Set-TextValue "D40" "0.0780"
$ws.Range("E40").Value = "  +9.99%  "

$ws.Range("E41").Value = "  +11.19%  "

Set-TextValue "D42" "2.926.68"
$ws.Range("E42").Value = "  +4.78%  "

$ws.Range("E43").Value = "  +2.33%  "

Set-TextValue "D44" "0.771"
$ws.Range("E44").Value = "  +7.20%  "

Set-TextValue "D45" "41.79"
$ws.Range("E45").Value = "  +4.80%  "

Set-TextValue "D46" "4.34"
$ws.Range("E46").Value = "  +3.01%  "

$ws.Range("E47").Value = "  +10.31%  "

Set-TextValue "D48" "3.487.07"
$ws.Range("E48").Value = "  +7.72%  "

$ws.Range("E49").Value = "  +9.38%  "

# Rows 50 and 51 swap Bittensor/Cosmos data along with updated price/volume
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "6.38"
$ws.Range("E50").Value = "  +3.58%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D51" "296.29"
$ws.Range("E51").Value = "  +12.41%  "
